# Update LR-pair sheet with new TPM-derived values.
# - Row 2 ("ECs" sending cluster) is removed entirely.
# - Old row 3 ("FAPs" sending cluster) becomes the new row 2 (specificity
#   columns I/J/S/T recompute because the total across rows changed).
# - Old row 4 ("Resolving-Mac" sending cluster) becomes the new row 3, with
#   its sending-cluster label renamed to "MuSCs" and fresh TPM-derived
#   expression/specificity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "ECs" row (current row 2); shift remaining rows up.
$ws.Range("A2:T2").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# New row 2 (previously "FAPs" row) — same raw values, recomputed specificity.
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl28"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1395456666666667
$ws.Range("H2").Value = 0.418637
$ws.Range("I2").Value = 0.5708284189068497
$ws.Range("J2").Value = 0.5708284189068498
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2746273333333333
$ws.Range("N2").Value = 0.823882
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.03832305431488889
$ws.Range("R2").Value = 0.344907488834
$ws.Range("S2").Value = 0.5708284189068497
$ws.Range("T2").Value = 0.5708284189068498

# New row 3 (previously "Resolving-Mac" row) — renamed sending cluster to
# "MuSCs" and updated TPM-derived expression/specificity values.
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Ccl28"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.104916
$ws.Range("H3").Value = 0.314748
$ws.Range("I3").Value = 0.4291715810931503
$ws.Range("J3").Value = 0.4291715810931503
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2746273333333333
$ws.Range("N3").Value = 0.823882
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.028812801304
$ws.Range("R3").Value = 0.259315211736
$ws.Range("S3").Value = 0.4291715810931503
$ws.Range("T3").Value = 0.4291715810931503
